# Applies the "deleting employees and adding employees" update to the
# George Smith 2020 leave registry sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("George Smith_2020_leave_registry")

# --- Header / summary fields -------------------------------------------------
$ws.Range("G5").Value = "August 24, 2000"
$ws.Range("G6").Value = 0
$ws.Range("G8").Value = "June 04, 2020"
$ws.Range("G10").Value = "June 04, 2020"

# --- Bi-monthly earned credits summary (rows 13-15) --------------------------
$ws.Range("H13").Value = 5
$ws.Range("I13").Value = 0.625

$ws.Range("H14").Value = 9
$ws.Range("I14").Value = -5.256

$ws.Range("C15").Value = 21
$ws.Range("D15").Value = 20

# --- Leave detail rows (19-21) ------------------------------------------------
$ws.Range("F19").Value = "May 26, 2020"
$ws.Range("G19").Value = 5
$ws.Range("I19").Value = "VL"

$ws.Range("F20").Value = "June 10, 2020"
$ws.Range("G20").Value = 4
$ws.Range("I20").Value = "SL"

$ws.Range("F21").Value = "June 11, 2020"
$ws.Range("G21").Value = 5
$ws.Range("I21").Value = "SL"

# --- Contract evaluation date -------------------------------------------------
$ws.Range("G43").Value = "August 24, 2019 to August 24, 2020"

# --- Offense log: update existing entry, then insert 3 new rows --------------
$ws.Range("F45").Value = "April 14, 2020"
$ws.Range("G45").Value = "Late"

$ws.Rows.Item(46).Resize(3).Insert()

$ws.Range("F46").Value = "June 10, 2020"
$ws.Range("G46").Value = "No Time-in"

$ws.Range("F47").Value = "June 16, 2020"
$ws.Range("G47").Value = "No Time-in"

$ws.Range("F48").Value = "June 09, 2020"
$ws.Range("G48").Value = "No Time-in"
